$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet, then refresh the Print_Area defined name so it points at the new sheet name
$ws.Name = "Step 255"
$ws.PageSetup.PrintArea = '$A$1:$B$25'

# Update formula in cell B23
$ws.Range("B23").Formula = "= (0.5)^3"

# Update cell A25 to be a formula that computes the chance string dynamically
$ws.Range("A25").Formula = "=""There is an approximately ""&ROUND(B25*100,2)&""% chance in binomial distribution that ""&B9&"" out of ""&B8&"" kids would prefer dogs over cats."""

$wb.Save()
